# Create Data Dictionary from DMN & BPMN models
# Populate the "Annotation" column (I) of the DMN data-dictionary table
# with per-row annotation notes, and clear the placeholder single-space
# value that previously sat in the table header's Annotation data cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 6 on the sheet / row 1 of the table) had a stray
# " " value in the Annotation column - clear it out entirely.
$ws.Range("I6").Value = ""

# Fill in annotation text for the data rows.
$ws.Range("I7").Value = "Annotation 02"
$ws.Range("I9").Value = "Annotation 04"
$ws.Range("I10").Value = "Annotation 05"
$ws.Range("I11").Value = "Annotation 06"

# Keep the active selection consistent with where editing finished.
$ws.Range("I8").Select()
